# Apply cryptocurrency price/volume updates per commit "Updated cryptos list" run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.412.28"
$ws.Range("E2").Value = "  +2.68%  "
$ws.Range("D3").Value = "2.067.11"
$ws.Range("E3").Value = "  +3.87%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'235.37"
$ws.Range("E5").Value = "  -0.48%  "
$ws.Range("D6").Value = "'0.614"
$ws.Range("E6").Value = "  +2.73%  "
$ws.Range("D7").Value = "'58.41"
$ws.Range("E7").Value = "  +7.19%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  +3.27%  "
$ws.Range("D10").Value = "'58.81"
$ws.Range("E10").Value = "  +1.09%  "
$ws.Range("E11").Value = "  +1.91%  "
$ws.Range("E12").Value = "  +3.46%  "
$ws.Range("D13").Value = "2.370.32"
$ws.Range("E13").Value = "  +3.86%  "
$ws.Range("E14").Value = "  +2.78%  "
$ws.Range("D15").Value = "'20.95"
$ws.Range("E15").Value = "  +3.68%  "
$ws.Range("D16").Value = "'0.780"
$ws.Range("E16").Value = "  +3.34%  "
$ws.Range("D17").Value = "'5.21"
$ws.Range("E17").Value = "  +2.86%  "
$ws.Range("D18").Value = "2.076.47"
$ws.Range("E18").Value = "  +5.82%  "
$ws.Range("D19").Value = "37.611.53"
$ws.Range("E19").Value = "  +3.20%  "
$ws.Range("D20").Value = "'6.19"
$ws.Range("E20").Value = "  +17.79%  "
$ws.Range("D21").Value = "'69.05"
$ws.Range("E21").Value = "  +1.77%  "
$ws.Range("D22").Value = "0.0₃0816"
$ws.Range("E22").Value = "  +1.66%  "
$ws.Range("D23").Value = "'226.66"
$ws.Range("E23").Value = "  +2.49%  "
$ws.Range("E24").Value = "  +0.15%  "
$ws.Range("E25").Value = "  +2.31%  "
$ws.Range("E26").Value = "  +1.14%  "
$ws.Range("D27").Value = "'164.60"
$ws.Range("E27").Value = "  +0.98%  "
$ws.Range("E28").Value = "  +14.10%  "
$ws.Range("E29").Value = "  +2.48%  "
$ws.Range("D30").Value = "'19.25"
$ws.Range("E30").Value = "  +1.82%  "
$ws.Range("E31").Value = "  -0.42%  "
$ws.Range("E32").Value = "  +2.30%  "
$ws.Range("E33").Value = "  +2.73%  "
$ws.Range("D34").Value = "'0.0623"
$ws.Range("E34").Value = "  +2.71%  "
$ws.Range("D35").Value = "'2.55"
$ws.Range("E35").Value = "  +9.01%  "
$ws.Range("E36").Value = "  +6.47%  "
$ws.Range("D37").Value = "'3.42"
$ws.Range("E37").Value = "  +3.69%  "
$ws.Range("E38").Value = "  +0.00%  "
$ws.Range("E39").Value = "  +0.57%  "
$ws.Range("E40").Value = "  +7.82%  "
$ws.Range("E41").Value = "  +7.63%  "
$ws.Range("B42").Value = "FTXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D42").Value = "'4.52"
$ws.Range("E42").Value = "  +24.64%  "
$ws.Range("B43").Value = "HuobiToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D43").Value = "'2.96"
$ws.Range("E43").Value = "  -1.05%  "
$ws.Range("D44").Value = "1.478.45"
$ws.Range("E44").Value = "  +1.65%  "
$ws.Range("D45").Value = "'96.78"
$ws.Range("E45").Value = "  +8.54%  "
$ws.Range("E46").Value = "  +6.01%  "
$ws.Range("D47").Value = "'0.0210"
$ws.Range("E47").Value = "  +4.37%  "
$ws.Range("D48").Value = "'15.90"
$ws.Range("E48").Value = "  +6.16%  "
$ws.Range("E49").Value = "  +3.65%  "
$ws.Range("E50").Value = "  +7.03%  "
$ws.Range("E51").Value = "  +1.70%  "
